# Simulator Overview - move persistence/publishing components to separate libraries.
# Slide 2 edits:
#  1. The second "CC2" label (between CC2 and CC4 in the calculation-context
#     timeline) is renamed to "CC3".
#  2. The "Market" label near the top-right legend is split into two runs
#     ("Market " + "Env") and the textbox is widened to fit the new text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- 1) CC2 -> CC3 -----------------------------------------------------
# (there are two "CC2" textboxes on this slide - "TextBox 92" and
# "TextBox 93"; only the second one, sitting between the CC2 and CC4
# markers on the timeline, is renamed)
$ccShape = $s.Shapes.Item("TextBox 93")
$ccShape.TextFrame.TextRange.Text = "CC3"

# --- 2) Market -> "Market " + "Env" (widen textbox) ---------------------
$marketShape = $s.Shapes.Item("TextBox 113")
$marketRange = $marketShape.TextFrame.TextRange
$marketRange.Text = "Market "
$marketRange.InsertAfter("Env") | Out-Null

# Widen the textbox (635751 EMU -> 928139 EMU; height unchanged) so the
# added "Env" run fits; COM Width/Height are expressed in points
# (928139 EMU / 12700 EMU-per-point ~= 73.0818110236pt; the literal below
# is nudged a hair above that exact quotient so the point -> EMU
# round-trip lands back on 928139 instead of 928138).
$marketShape.Width = 73.08185
